$d = $word.ActiveDocument

# --- Table 3: "PERSONA QUE OCUPA DE MANERA PROVISIONAL LA PLAZA" ---
# Header row has three merged "EJEMPLO" placeholder cells sitting above
# APELLIDO PATERNO / APELLIDO MATERNO / NOMBRE(S); clear them to blank.
# (Each placeholder cell spans several grid columns because of w:gridSpan,
# so Cell(1,1)/Cell(1,3)/Cell(1,6) land on the three distinct <w:tc>s.)
$tEjemplo = $d.Tables.Item(3)
$tEjemplo.Cell(1, 1).Range.Text = ""
$tEjemplo.Cell(1, 3).Range.Text = ""
$tEjemplo.Cell(1, 6).Range.Text = ""

# --- Table 4: "A PARTIR DEL:" (reinstatement date) ---
# DÍA / MES / AÑO values: 02/FEBRERO/2026 -> 19/NOVIEMBRE/2025
$tFecha = $d.Tables.Item(4)
$tFecha.Cell(1, 2).Range.Text = "19"
$tFecha.Cell(1, 4).Range.Text = "NOVIEMBRE"
$tFecha.Cell(1, 6).Range.Text = "2025"
